## Fruta / hortaliza, semanal
## Insert two new weekly price-report rows (Primera/Segunda) for
## Pepino ensalada before the existing row 340, pushing the rest of the
## table down by two rows (old row 340 -> new row 342, etc.). The table's
## previous last two rows (old 440/441) end up as the new last two rows
## (442/443).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 340:341 - this shifts existing rows 340-441
# down to 342-443, carrying their formatting (e.g. the date style on
# column D) along with them.
$ws.Rows("340:341").Insert()

# New row 340 (Primera)
$ws.Range("A340").Value = 1
$ws.Range("B340").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C340").Value = "Arica y Parinacota"
$ws.Range("D340").Value = 44985
$ws.Range("E340").Value = 15
$ws.Range("F340").Value = 100112043
$ws.Range("G340").Value = "Pepino ensalada"
$ws.Range("H340").Value = "Sin especificar"
$ws.Range("I340").Value = "Primera"
$ws.Range("J340").Value = 300
$ws.Range("K340").Value = 3000
$ws.Range("L340").Value = 3500
$ws.Range("M340").Value = 3333
$ws.Range("N340").Value = "$/caja 70 unidades"
$ws.Range("O340").Value = "Región de Arica y Parinacota"
$ws.Range("P340").Value = 48
$ws.Range("Q340").Value = 70
$ws.Range("R340").Value = "Hortaliza"

# New row 341 (Segunda)
$ws.Range("A341").Value = 1
$ws.Range("B341").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C341").Value = "Arica y Parinacota"
$ws.Range("D341").Value = 44985
$ws.Range("E341").Value = 15
$ws.Range("F341").Value = 100112043
$ws.Range("G341").Value = "Pepino ensalada"
$ws.Range("H341").Value = "Sin especificar"
$ws.Range("I341").Value = "Segunda"
$ws.Range("J341").Value = 150
$ws.Range("K341").Value = 2500
$ws.Range("L341").Value = 3000
$ws.Range("M341").Value = 2800
$ws.Range("N341").Value = "$/caja 100 unidades"
$ws.Range("O341").Value = "Región de Arica y Parinacota"
$ws.Range("P341").Value = 28
$ws.Range("Q341").Value = 100
$ws.Range("R341").Value = "Hortaliza"
